# corrected again capacity market
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "times": shift the simulation StartTime/StopTime back to 2021/2022
# ---------------------------------------------------------------
$wsTimes = $wb.Worksheets.Item("times")
$wsTimes.Cells.Item(2, 2).Value = 44196.99861111111
$wsTimes.Cells.Item(3, 2).Value = 44560.99861111111

# ---------------------------------------------------------------
# Sheet "scenario_data_emlab": add a second year column (2020/2021),
# update fuel/CO2 prices for the (now first) year column.
# ---------------------------------------------------------------
$wsScenario = $wb.Worksheets.Item("scenario_data_emlab")

# Year headers
$wsScenario.Cells.Item(1, 2).Value = 2020
$wsScenario.Cells.Item(1, 3).Value = 2021

# Updated values in column B
$wsScenario.Cells.Item(2, 2).Value = 24.38
$wsScenario.Cells.Item(5, 2).Value = 11.504
$wsScenario.Cells.Item(6, 2).Value = 20.468
$wsScenario.Cells.Item(7, 2).Value = 48.114

# New column C: empty (but typed-as-text) placeholder cells for rows 2-7
$wsScenario.Cells.Item(2, 3).Value = "'"
$wsScenario.Cells.Item(3, 3).Value = "'"
$wsScenario.Cells.Item(4, 3).Value = "'"
$wsScenario.Cells.Item(5, 3).Value = "'"
$wsScenario.Cells.Item(6, 3).Value = "'"
$wsScenario.Cells.Item(7, 3).Value = "'"

# Row 8 demand series path repeated in the new column
$wsScenario.Cells.Item(8, 3).Value = "./timeseries/demand/load.csv"

# ---------------------------------------------------------------
# Sheet "conventionals": drop the two placeholder plants (rows 2-3)
# ---------------------------------------------------------------
$wsConv = $wb.Worksheets.Item("conventionals")

$wsConv.Cells.Item(2, 1).Value = 0
$wsConv.Cells.Item(2, 2).Value = 20150300022
$wsConv.Cells.Item(2, 3).Value = "NATURAL_GAS"
$wsConv.Cells.Item(2, 4).Value = 4.2
$wsConv.Cells.Item(2, 5).Value = 0.61
$wsConv.Cells.Item(2, 6).Value = 31358.329
$wsConv.Cells.Item(2, 7).Value = 31358.329

$wsConv.Cells.Item(3, 1).Value = 1
$wsConv.Cells.Item(3, 2).Value = 20152800024
$wsConv.Cells.Item(3, 3).Value = "HARD_COAL"
$wsConv.Cells.Item(3, 4).Value = 3.5
$wsConv.Cells.Item(3, 5).Value = 0.33
$wsConv.Cells.Item(3, 6).Value = 24845.77
$wsConv.Cells.Item(3, 7).Value = 24845.77

$wsConv.Cells.Item(4, 1).Value = 2
$wsConv.Cells.Item(4, 2).Value = 20153000025
$wsConv.Cells.Item(4, 3).Value = "OIL"
$wsConv.Cells.Item(4, 4).Value = 6
$wsConv.Cells.Item(4, 5).Value = 0.35
$wsConv.Cells.Item(4, 6).Value = 3652.9
$wsConv.Cells.Item(4, 7).Value = 3652.9

$wsConv.Cells.Item(5, 1).Value = 3
$wsConv.Cells.Item(5, 2).Value = 20152900027
$wsConv.Cells.Item(5, 3).Value = "LIGNITE"
$wsConv.Cells.Item(5, 4).Value = 3.5
$wsConv.Cells.Item(5, 5).Value = 0.33
$wsConv.Cells.Item(5, 6).Value = 20779.02
$wsConv.Cells.Item(5, 7).Value = 20779.02

$wsConv.Cells.Item(6, 1).Value = 4
$wsConv.Cells.Item(6, 2).Value = 20151400028
$wsConv.Cells.Item(6, 3).Value = "NUCLEAR"
$wsConv.Cells.Item(6, 4).Value = 3.5
$wsConv.Cells.Item(6, 5).Value = 0.33
$wsConv.Cells.Item(6, 6).Value = 8599
$wsConv.Cells.Item(6, 7).Value = 8599

$wsConv.Cells.Item(7, 1).Value = 5
$wsConv.Cells.Item(7, 2).Value = 20151700029
$wsConv.Cells.Item(7, 3).Value = "NATURAL_GAS"
$wsConv.Cells.Item(7, 4).Value = 4.5
$wsConv.Cells.Item(7, 5).Value = 0.43
$wsConv.Cells.Item(7, 6).Value = 8194.3025
$wsConv.Cells.Item(7, 7).Value = 8194.3025

# rows 8 & 9 (old NUCLEAR / NATURAL_GAS rows, now duplicated above) go away
$wsConv.Rows.Item(9).Delete()
$wsConv.Rows.Item(8).Delete()

# ---------------------------------------------------------------
# Sheet "renewables": replace the two PV/Wind placeholders with a
# real WindOn plant, and add nine new 1000 MW OtherPV capacity-market
# entries ahead of the existing RunOfRiver/PV/WindOff plants.
# ---------------------------------------------------------------
$wsRen = $wb.Worksheets.Item("renewables")

# row 2 (index 0) becomes the real WindOn plant (former row 5)
$wsRen.Cells.Item(2, 2).Value = 20152400023
$wsRen.Cells.Item(2, 3).Value = 47547.50848700004
$wsRen.Cells.Item(2, 4).Value = 1.35
$wsRen.Cells.Item(2, 5).Value = "WindOn"

# row 3 (index 1)
$wsRen.Cells.Item(3, 2).Value = 20212100031
$wsRen.Cells.Item(3, 3).Value = 1000
$wsRen.Cells.Item(3, 4).Value = 0
$wsRen.Cells.Item(3, 5).Value = "OtherPV"

# row 4 (index 2)
$wsRen.Cells.Item(4, 2).Value = 20212100032
$wsRen.Cells.Item(4, 3).Value = 1000
$wsRen.Cells.Item(4, 4).Value = 0
$wsRen.Cells.Item(4, 5).Value = "OtherPV"

# row 5 (index 3)
$wsRen.Cells.Item(5, 2).Value = 20212100034
$wsRen.Cells.Item(5, 3).Value = 1000
$wsRen.Cells.Item(5, 4).Value = 0
$wsRen.Cells.Item(5, 5).Value = "OtherPV"

# row 6 (index 4)
$wsRen.Cells.Item(6, 2).Value = 20212100037
$wsRen.Cells.Item(6, 3).Value = 1000
$wsRen.Cells.Item(6, 4).Value = 0
$wsRen.Cells.Item(6, 5).Value = "OtherPV"

# row 7 (index 5)
$wsRen.Cells.Item(7, 2).Value = 20212100041
$wsRen.Cells.Item(7, 3).Value = 1000
$wsRen.Cells.Item(7, 4).Value = 0
$wsRen.Cells.Item(7, 5).Value = "OtherPV"

# insert 7 fresh rows (8-14) before the old row 8 (RunOfRiver/PV/WindOff block),
# cloning formatting from row 7 so borders/styles stay consistent
$wsRen.Rows.Item(7).Copy()
$wsRen.Rows.Item(8).Insert()
$wsRen.Rows.Item(7).Copy()
$wsRen.Rows.Item(8).Insert()
$wsRen.Rows.Item(7).Copy()
$wsRen.Rows.Item(8).Insert()
$wsRen.Rows.Item(7).Copy()
$wsRen.Rows.Item(8).Insert()
$wsRen.Rows.Item(7).Copy()
$wsRen.Rows.Item(8).Insert()
$wsRen.Rows.Item(7).Copy()
$wsRen.Rows.Item(8).Insert()
$wsRen.Rows.Item(7).Copy()
$wsRen.Rows.Item(8).Insert()

# row 8 (index 6)
$wsRen.Cells.Item(8, 1).Value = 6
$wsRen.Cells.Item(8, 2).Value = 20212100046
$wsRen.Cells.Item(8, 3).Value = 1000
$wsRen.Cells.Item(8, 4).Value = 0
$wsRen.Cells.Item(8, 5).Value = "OtherPV"
$wsRen.Cells.Item(8, 6).Value = "-"
$wsRen.Cells.Item(8, 7).Value = "-"
$wsRen.Cells.Item(8, 8).Value = "-"
$wsRen.Cells.Item(8, 9).Value = "-"

# row 9 (index 7)
$wsRen.Cells.Item(9, 1).Value = 7
$wsRen.Cells.Item(9, 2).Value = 20212100052
$wsRen.Cells.Item(9, 3).Value = 1000
$wsRen.Cells.Item(9, 4).Value = 0
$wsRen.Cells.Item(9, 5).Value = "OtherPV"
$wsRen.Cells.Item(9, 6).Value = "-"
$wsRen.Cells.Item(9, 7).Value = "-"
$wsRen.Cells.Item(9, 8).Value = "-"
$wsRen.Cells.Item(9, 9).Value = "-"

# row 10 (index 8)
$wsRen.Cells.Item(10, 1).Value = 8
$wsRen.Cells.Item(10, 2).Value = 20212100059
$wsRen.Cells.Item(10, 3).Value = 1000
$wsRen.Cells.Item(10, 4).Value = 0
$wsRen.Cells.Item(10, 5).Value = "OtherPV"
$wsRen.Cells.Item(10, 6).Value = "-"
$wsRen.Cells.Item(10, 7).Value = "-"
$wsRen.Cells.Item(10, 8).Value = "-"
$wsRen.Cells.Item(10, 9).Value = "-"

# row 11 (index 9)
$wsRen.Cells.Item(11, 1).Value = 9
$wsRen.Cells.Item(11, 2).Value = 20212100067
$wsRen.Cells.Item(11, 3).Value = 1000
$wsRen.Cells.Item(11, 4).Value = 0
$wsRen.Cells.Item(11, 5).Value = "OtherPV"
$wsRen.Cells.Item(11, 6).Value = "-"
$wsRen.Cells.Item(11, 7).Value = "-"
$wsRen.Cells.Item(11, 8).Value = "-"
$wsRen.Cells.Item(11, 9).Value = "-"

# row 12 (index 10)
$wsRen.Cells.Item(12, 1).Value = 10
$wsRen.Cells.Item(12, 2).Value = 20212100076
$wsRen.Cells.Item(12, 3).Value = 1000
$wsRen.Cells.Item(12, 4).Value = 0
$wsRen.Cells.Item(12, 5).Value = "OtherPV"
$wsRen.Cells.Item(12, 6).Value = "-"
$wsRen.Cells.Item(12, 7).Value = "-"
$wsRen.Cells.Item(12, 8).Value = "-"
$wsRen.Cells.Item(12, 9).Value = "-"

# row 13 (index 11) - former RunOfRiver plant
$wsRen.Cells.Item(13, 1).Value = 11
$wsRen.Cells.Item(13, 2).Value = 20151200026
$wsRen.Cells.Item(13, 3).Value = 8858.749999999998
$wsRen.Cells.Item(13, 4).Value = 0
$wsRen.Cells.Item(13, 5).Value = "RunOfRiver"
$wsRen.Cells.Item(13, 6).Value = "-"
$wsRen.Cells.Item(13, 7).Value = "-"
$wsRen.Cells.Item(13, 8).Value = "-"
$wsRen.Cells.Item(13, 9).Value = "-"

# row 14 (index 12) - former large OtherPV plant
$wsRen.Cells.Item(14, 1).Value = 12
$wsRen.Cells.Item(14, 2).Value = 20152100030
$wsRen.Cells.Item(14, 3).Value = 53555.51607579708
$wsRen.Cells.Item(14, 4).Value = 0
$wsRen.Cells.Item(14, 5).Value = "OtherPV"
$wsRen.Cells.Item(14, 6).Value = "-"
$wsRen.Cells.Item(14, 7).Value = "-"
$wsRen.Cells.Item(14, 8).Value = "-"
$wsRen.Cells.Item(14, 9).Value = "-"

# row 15 (index 13) is the former row 8 (WindOff) shifted down by the inserts;
# only its index needs correcting, the rest of the row is untouched
$wsRen.Cells.Item(15, 1).Value = 13

# ---------------------------------------------------------------
# Sheet "storages": drop the single placeholder STORAGE row
# ---------------------------------------------------------------
$wsStor = $wb.Worksheets.Item("storages")
$wsStor.Rows.Item(2).Delete()

# ---------------------------------------------------------------
# Sheet "biogas": drop the placeholder plant, keep the real one
# ---------------------------------------------------------------
$wsBiogas = $wb.Worksheets.Item("biogas")
$wsBiogas.Cells.Item(2, 2).Value = 20150100021
$wsBiogas.Cells.Item(2, 3).Value = 4644.4034
$wsBiogas.Rows.Item(3).Delete()

Write-Host "edit applied"
